# LoginPage Pom class added
# Replace the old sample data (Akshara / Swara) with a username/password
# login-form style layout and widen column A slightly, matching the
# "LoginPage" Page-Object-Model fixture data used by the Selenium tests.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: existing "sunday" label stays in A1, new "password" header in B1
$ws.Range("A1").Value = "sunday"
$ws.Range("B1").Value = "password"

# Row 2: replace "Swara" with the Admin/admin123 login credentials
$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "admin123"

# Column A gets an explicit width (~10 characters)
$ws.Columns.Item(1).ColumnWidth = 9.14

# Leave the selection on B2, the last-edited cell
$ws.Range("B2").Select() | Out-Null
